{"js": "const replacements = [\n  [\"2024-05-17 Friday\", \"2024-05-18 Saturday\"],\n  [\"697\u00f75=\", \"924\u00f72=\"],\n  [\"520\u00f72=\", \"889\u00f72=\"],\n  [\"336\u00f78=\", \"370\u00f72=\"],\n  [\"778\u00f77=\", \"977\u00f77=\"],\n  [\"318\u00f76=\", \"470\u00f76=\"],\n  [\"854\u00f79=\", \"460\u00f73=\"],\n  [\"505\u00f76=\", \"111\u00f78=\"],\n  [\"133\u00f77=\", \"836\u00f72=\"],\n  [\"119\u00f75=\", \"409\u00f78=\"],\n  [\"699\u00f76=\", \"452\u00f78=\"],\n  [\"777\u00f77=\", \"351\u00f74=\"],\n  [\"140\u00f73=\", \"741\u00f76=\"],\n  [\"595\u00f77=\", \"381\u00f79=\"],\n  [\"224\u00f74=\", \"296\u00f79=\"],\n  [\"265\u00f75=\", \"669\u00f72=\"],\n  [\"623\u00f73=\", \"531\u00f72=\"],\n  [\"829\u00f74=\", \"335\u00f72=\"],\n  [\"145\u00f76=\", \"815\u00f78=\"],\n  [\"857\u00f79=\", \"486\u00f72=\"],\n  [\"831\u00f79=\", \"865\u00f78=\"],\n  [\"586\u00f78=\", \"381\u00f76=\"],\n  [\"772\u00f78=\", \"507\u00f72=\"],\n  [\"374\u00f72=\", \"679\u00f78=\"],\n  [\"790\u00f79=\", \"472\u00f75=\"],\n  [\"105\u00f78=\", \"676\u00f73=\"],\n];\n\nfor (const [from, to] of replacements) {\n  const results = context.document.body.search(from, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(to, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$replacements = @(\n    @(\"2024-05-17 Friday\", \"2024-05-18 Saturday\"),\n    @(\"697\u00f75=\", \"924\u00f72=\"),\n    @(\"520\u00f72=\", \"889\u00f72=\"),\n    @(\"336\u00f78=\", \"370\u00f72=\"),\n    @(\"778\u00f77=\", \"977\u00f77=\"),\n    @(\"318\u00f76=\", \"470\u00f76=\"),\n    @(\"854\u00f79=\", \"460\u00f73=\"),\n    @(\"505\u00f76=\", \"111\u00f78=\"),\n    @(\"133\u00f77=\", \"836\u00f72=\"),\n    @(\"119\u00f75=\", \"409\u00f78=\"),\n    @(\"699\u00f76=\", \"452\u00f78=\"),\n    @(\"777\u00f77=\", \"351\u00f74=\"),\n    @(\"140\u00f73=\", \"741\u00f76=\"),\n    @(\"595\u00f77=\", \"381\u00f79=\"),\n    @(\"224\u00f74=\", \"296\u00f79=\"),\n    @(\"265\u00f75=\", \"669\u00f72=\"),\n    @(\"623\u00f73=\", \"531\u00f72=\"),\n    @(\"829\u00f74=\", \"335\u00f72=\"),\n    @(\"145\u00f76=\", \"815\u00f78=\"),\n    @(\"857\u00f79=\", \"486\u00f72=\"),\n    @(\"831\u00f79=\", \"865\u00f78=\"),\n    @(\"586\u00f78=\", \"381\u00f76=\"),\n    @(\"772\u00f78=\", \"507\u00f72=\"),\n    @(\"374\u00f72=\", \"679\u00f78=\"),\n    @(\"790\u00f79=\", \"472\u00f75=\"),\n    @(\"105\u00f78=\", \"676\u00f73=\"),\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $from = $pair[0]\n    $to = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $from\n    $find.Replacement.Text = $to\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
